# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header values -------------------------------------------------------
# VALOR MORA total grows with the new period added
$ws.Range("E11").Value = 520000
# Cant. Periodos grows from 9 to 10
$ws.Range("F13").Value = 10

# Swap the "Novedad de Ingreso" / "Novedad de Retiro" column headers
$ws.Range("H15").Value = "Novedad de Retiro"
$ws.Range("I15").Value = "Novedad de Ingreso"

# --- Insert a new period row (2509) ---------------------------------------
# Row 24 (period 2508) is currently the last data row and carries the
# special "closing" bottom border. Insert a fresh row below it, copy row
# 24's format+values down into it, then restyle row 24 like the other
# interior rows (16-23) and set the new period value.
$ws.Rows("25:25").Insert()
$ws.Range("B24:J24").Copy($ws.Range("B25:J25"))
$ws.Range("E25").Value = "2509"

$ws.Range("B23:J23").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Center-align the "Periodo Mora" column for every data row (16-25)
$ws.Range("E16:E25").HorizontalAlignment = -4108
